$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seguimiento")

# Row 8 (guion 3) first, so the new shared-string entries land in the same
# order as the target file: index 14 = "Manuscritos enviados..." (F8),
# index 15 = "Archivos completos..." (F7).
$ws.Range("B8").Value = 42080
$ws.Range("C8").Value = 42080
$ws.Range("F8").Value = "Manuscritos enviados 16/03/2015 a corrección estilo"

# Row 7 (guion 2): manuscript/format/skeleton/graphic-request dates + "entrega" note
$ws.Range("B7").Value = 42079
$ws.Range("C7").Value = 42079
$ws.Range("D7").Value = 42079
$ws.Range("E7").Value = 42079
$ws.Range("F7").Value = "Archivos completos con corrección estilo"

# The longer wrapped "entrega" notes need more vertical room than the
# default row height, same as Excel grows the row automatically when the
# wrapped text no longer fits on one line.
$ws.Rows.Item(7).RowHeight = 30.75
$ws.Rows.Item(8).RowHeight = 45.75

# The rest of the data rows share the sheet's normal (non-wrapped) row
# height; align them on the same 16.5pt baseline as the edited rows.
foreach ($r in 2,3,6,9,10,11,12,13) {
    $ws.Rows.Item($r).RowHeight = 16.5
}

# Update the active view selection/scroll position as recorded after the edit
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("E9").Select()
